$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '67.053.54'
$ws.Range("E2").Value = '  -1.30%  '

# Row 3
$ws.Range("D3").Value = '3.518.55'

# Row 4
$style = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = $style
$ws.Range("E4").Value = '  -0.05%  '

# Row 5
$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '609.73'
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = '  +0.48%  '

# Row 6
$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.22'
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = '  -1.69%  '

# Row 7
$ws.Range("D7").Value = '3.517.06'
$ws.Range("E7").Value = '  +0.40%  '

# Row 8
$ws.Range("E8").Value = '  +0.03%  '

# Row 9
$ws.Range("E9").Value = '  -1.60%  '

# Row 10
$ws.Range("E10").Value = '  -1.02%  '

# Row 11
$style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '8.07'
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = '  +6.47%  '

# Row 12
$style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.424'
$ws.Range("D12").Style = $style
$ws.Range("E12").Value = '  -1.67%  '

# Row 14
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '4.113.58'
$ws.Range("E14").Value = '  +0.43%  '

# Row 15
$ws.Range("B15").Value = 'Avalanche'
$ws.Range("C15").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '31.99'
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = '  -0.20%  '

# Row 16
$ws.Range("D16").Value = '3.518.82'
$ws.Range("E16").Value = '  +0.60%  '

# Row 17
$ws.Range("D17").Value = '67.084.16'
$ws.Range("E17").Value = '  -1.24%  '

# Row 18
$ws.Range("E18").Value = '  -0.16%  '

# Row 19
$style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.84'
$ws.Range("D19").Style = $style
$ws.Range("E19").Value = '  +8.88%  '

# Row 20
$style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.48'
$ws.Range("D20").Style = $style
$ws.Range("E20").Value = '  -0.07%  '

# Row 21
$style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.43'
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = '  +0.21%  '

# Row 22
$style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '438.53'
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = '  -1.78%  '

# Row 23
$style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.611'
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = '  -2.34%  '

# Row 24
$style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '79.82'
$ws.Range("D24").Style = $style
$ws.Range("E24").Value = '  +0.91%  '

# Row 25
$ws.Range("D25").Value = '3.652.61'
$ws.Range("E25").Value = '  +0.24%  '

# Row 26
$ws.Range("E26").Value = '  -0.02%  '

# Row 27
$ws.Range("E27").Value = '  -3.92%  '

# Row 28
$style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.82'
$ws.Range("D28").Style = $style
$ws.Range("E28").Value = '  -1.77%  '

# Row 29
$style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.31'
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = '  -3.57%  '

# Row 30
$ws.Range("E30").Value = '  +0.61%  '

# Row 31
$ws.Range("E31").Value = '  -1.97%  '

# Row 32
$ws.Range("E32").Value = '  -1.53%  '

# Row 33
$ws.Range("E33").Value = '  -0.04%  '

# Row 34
$style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '25.65'
$ws.Range("D34").Style = $style
$ws.Range("E34").Value = '  +0.10%  '

# Row 35
$style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.99'
$ws.Range("D35").Style = $style
$ws.Range("E35").Value = '  -2.64%  '

# Row 36
$ws.Range("E36").Value = '  -2.03%  '

# Row 37
$ws.Range("E37").Value = '  +1.34%  '

# Row 39
$style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.998'
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = '  -0.04%  '

# Row 40
$style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '175.87'
$ws.Range("D40").Style = $style
$ws.Range("E40").Value = '  -0.49%  '

# Row 41
$style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0898'
$ws.Range("D41").Style = $style
$ws.Range("E41").Value = '  -0.13%  '

# Row 42
$ws.Range("E42").Value = '  -0.10%  '

# Row 43
$style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.07'
$ws.Range("D43").Style = $style
$ws.Range("E43").Value = '  -11.52%  '

# Row 44
$ws.Range("E44").Value = '  -0.09%  '

# Row 45
$style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '46.20'
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = '  -1.07%  '

# Row 46
$style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '28.20'
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = '  -7.49%  '

# Row 47
$ws.Range("E47").Value = '  -2.41%  '

# Row 48
$ws.Range("B48").Value = 'dogwifhat'
$ws.Range("C48").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.47'
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = '  -2.59%  '

# Row 49
$ws.Range("B49").Value = 'Cosmos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.49'
$ws.Range("D49").Style = $style
$ws.Range("E49").Value = '  -1.68%  '

# Row 50
$style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.999'
$ws.Range("D50").Style = $style
$ws.Range("E50").Value = '  +0.81%  '

# Row 51
$ws.Range("E51").Value = '  -1.53%  '
